$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 28: 2018-11-15 (serial 43419), 1 hour, "Meeting"
$ws.Range("A28").Value = 43419
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = "Meeting"

# New row 29: 2018-11-17 (serial 43421), 2 hours, new description
$ws.Range("A29").Value = 43421
$ws.Range("A27").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = "In order brengen fietssimulatie + classificatie toevoegen"

# Update selection to match target state
$ws.Range("B30").Select()
